# Commit: "Classify specific pairs of code sections within case as superstrike [#166297321]"
#
# Adds a new case record ("LECTER, HANNIBAL", CII/SSN 1009123456) with three conviction
# rows (35-37) to the gogen_pilots Los Angeles test fixture. The second conviction row
# (36) carries the case-level "Superstrike Code Section(s)" classification for the
# paired code sections "136.1 PC + 186.22(B)(4) PC".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new rows with the formatting of the most similar existing case-group rows
# (row 25 for the first conviction in a group, row 26 for subsequent ones) before
# filling in the data, matching how this sheet's row-groups are normally built.
$ws.Range("A25:DH25").Copy()
$ws.Range("A35:DH35").PasteSpecial(-4122)
$ws.Range("A26:DH26").Copy()
$ws.Range("A36:DH36").PasteSpecial(-4122)
$ws.Range("A26:DH26").Copy()
$ws.Range("A37:DH37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 35
$ws.Range("D35").Value = 84734893
$ws.Range("E35").Value = "'" + '#'
$ws.Range("F35").Value = "'" + '1009123456'
$ws.Range("M35").Value = "'" + '1009123456'
$ws.Range("N35").Value = "'" + 'LECTER, HANNIBAL'
$ws.Range("P35").Value = 19721127
$ws.Range("Q35").Value = 859349027
$ws.Range("R35").Value = "'" + 'x       '
$ws.Range("AM35").Value = "'" + 'x       '
$ws.Range("AP35").Value = 20150214
$ws.Range("AR35").Value = "'" + 'COURT ACTION'
$ws.Range("AX35").Value = "'" + 'LOS ANGELES'
$ws.Range("AZ35").Value = 101001022000
$ws.Range("BA35").Value = 20150315
$ws.Range("BC35").Value = "'" + 'x'
$ws.Range("BD35").Value = "'" + '136.1 PC-SPYING ON CATS'
$ws.Range("BE35").Value = "'" + 'F'
$ws.Range("CD35").Value = "'" + 'CONVICTED-PROBATION'
$ws.Range("CF35").Value = "'" + 'FELONY'
$ws.Range("CG35").Value = "'" + '#'
$ws.Range("CJ35").Value = "'" + 'PROBATION'
$ws.Range("CK35").Value = 2
$ws.Range("CL35").Value = "'" + 'M'
$ws.Range("CM35").Value = "'" + 'MONTHS'
$ws.Range("CN35").Value = 42

# Row 36
$ws.Range("A36").Value = "'" + 'true'
$ws.Range("D36").Value = 84734893
$ws.Range("E36").Value = "'" + '#'
$ws.Range("F36").Value = "'" + '1009123456'
$ws.Range("M36").Value = "'" + '1009123456'
$ws.Range("N36").Value = "'" + 'LECTER, HANNIBAL'
$ws.Range("P36").Value = 19721127
$ws.Range("Q36").Value = 859349027
$ws.Range("R36").Value = "'" + 'x       '
$ws.Range("AM36").Value = "'" + 'x       '
$ws.Range("AP36").Value = 20150519
$ws.Range("AR36").Value = "'" + 'COURT ACTION'
$ws.Range("AX36").Value = "'" + 'LOS ANGELES'
$ws.Range("AZ36").Value = 101001023000
$ws.Range("BA36").Value = 20150522
$ws.Range("BC36").Value = "'" + 'x'
$ws.Range("BD36").Value = "'" + '11358 HS-CULTIVATE CANNABIS'
$ws.Range("BE36").Value = "'" + 'F'
$ws.Range("CD36").Value = "'" + 'CONVICTED-PROBATION'
$ws.Range("CF36").Value = "'" + 'FELONY'
$ws.Range("CG36").Value = "'" + '#'
$ws.Range("CI36").Value = "'" + 'P'
$ws.Range("CJ36").Value = "'" + 'PROBATION'
$ws.Range("CK36").Value = 2
$ws.Range("CL36").Value = "'" + 'M'
$ws.Range("CM36").Value = "'" + 'MONTHS'
$ws.Range("CN36").Value = 42
$ws.Range("CT36").Value = "'" + '3'
$ws.Range("CU36").Value = "'" + '136.1 PC + 186.22(B)(4) PC'
$ws.Range("CV36").Value = "'" + '-'
$ws.Range("CW36").Value = "'" + '-'
$ws.Range("CX36").Value = "'" + '05/19/2015'
$ws.Range("CY36").Value = "'" + '4.5'
$ws.Range("CZ36").Value = "'" + '4.0'
$ws.Range("DA36").Value = "'" + '1'
$ws.Range("DB36").Value = "'" + '0'
$ws.Range("DC36").Value = "'" + '1'
$ws.Range("DD36").Value = "'" + '0'
$ws.Range("DE36").Value = "'" + '0'
$ws.Range("DF36").Value = "'" + '-'
$ws.Range("DG36").Value = "'" + 'Not eligible'
$ws.Range("DH36").Value = "'" + 'PC 667(e)(2)(c)(iv)'

# Row 37
$ws.Range("D37").Value = 84734893
$ws.Range("E37").Value = "'" + '#'
$ws.Range("F37").Value = "'" + '1009123456'
$ws.Range("M37").Value = "'" + '1009123456'
$ws.Range("N37").Value = "'" + 'LECTER, HANNIBAL'
$ws.Range("P37").Value = 19721127
$ws.Range("Q37").Value = 859349027
$ws.Range("R37").Value = "'" + 'x       '
$ws.Range("AM37").Value = "'" + 'x       '
$ws.Range("AP37").Value = 20151031
$ws.Range("AR37").Value = "'" + 'COURT ACTION'
$ws.Range("AX37").Value = "'" + 'LOS ANGELES'
$ws.Range("AZ37").Value = 101001024000
$ws.Range("BA37").Value = 20151031
$ws.Range("BC37").Value = "'" + 'x'
$ws.Range("BD37").Value = "'" + '186.22(B)(4) PC- GANG STUFF'
$ws.Range("BE37").Value = "'" + 'F'
$ws.Range("CD37").Value = "'" + 'CONVICTED-PROBATION'
$ws.Range("CF37").Value = "'" + 'FELONY'
$ws.Range("CG37").Value = "'" + '#'
$ws.Range("CI37").Value = "'" + 'P'
$ws.Range("CJ37").Value = "'" + 'PROBATION'
$ws.Range("CK37").Value = 2
$ws.Range("CL37").Value = "'" + 'M'
$ws.Range("CM37").Value = "'" + 'MONTHS'
$ws.Range("CN37").Value = 42

# Leave the selection where the editor ended up after data entry
$ws.Range("BE48").Select()
